# LinkedIn_Analysis_Report.xlsx — add two new PMPs (Jacqueline Shen, Nikki Gittins)
# discovered by the flexible all-PMP assignment mode, and refresh the Summary
# sheet's aggregate metrics to reflect the now-larger PMP roster.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: URL_Validation
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("URL_Validation")

$jacURL = "https://www.linkedin.com/in/jacqueline-shen-5b565917a/"

# Row 22 - Jacqueline Shen (valid LinkedIn URL)
$ws1.Range("A22").Value = "Jacqueline Shen"

$ws1.Range("B22").Value = $jacURL
$ws1.Hyperlinks.Add($ws1.Range("B22"), $jacURL)
$ws1.Range("B22").Style = $ws1.Range("B3").Style()

$ws1.Range("C22").Value = $true

$ws1.Range("D22").Value = $jacURL
$ws1.Hyperlinks.Add($ws1.Range("D22"), $jacURL)
$ws1.Range("D22").Style = $ws1.Range("D3").Style()

$ws1.Range("E22").Value = 10
$ws1.Range("F22").Value = "[]"

# Row 23 - Nikki Gittins (no LinkedIn URL provided)
$ws1.Range("A23").Value = "Nikki Gittins"
$ws1.Range("B23").Value = "nan"
$ws1.Range("C23").Value = $false
$ws1.Range("E23").Value = 0
$ws1.Range("F23").Value = "['No LinkedIn URL provided']"

# ---------------------------------------------------------------------------
# Sheet 2: LinkedIn_Analysis
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LinkedIn_Analysis")

# Row 22 - Jacqueline Shen
$ws2.Range("A22").Value = "Jacqueline Shen"

$ws2.Range("B22").Value = $jacURL
$ws2.Hyperlinks.Add($ws2.Range("B22"), $jacURL)
$ws2.Range("B22").Style = $ws2.Range("B3").Style()

$ws2.Range("C22").Value = 10
$ws2.Range("D22").Value = 9
$ws2.Range("E22").Value = 4.37
$ws2.Range("F22").Value = "Yes"
$ws2.Range("G22").Value = "High"

# Row 23 - Nikki Gittins
$ws2.Range("A23").Value = "Nikki Gittins"
$ws2.Range("C23").Value = 0
$ws2.Range("D23").Value = 8
$ws2.Range("E23").Value = 4.94
$ws2.Range("F23").Value = "No"
$ws2.Range("G23").Value = "High"

# ---------------------------------------------------------------------------
# Sheet 3: Summary — refresh aggregate metrics for the 22-PMP roster
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Summary")

$ws3.Range("B2").Value = 22   # Total PMPs
$ws3.Range("B3").Value = 18   # Valid LinkedIn URLs

# These metrics are stored as text in the workbook (not numbers), so force
# the Text number format before assigning to avoid Excel's automatic
# number/percentage coercion on entry.
$ws3.Range("B4").NumberFormat = "@"
$ws3.Range("B4").Value = "81.8%"   # LinkedIn Coverage %

$ws3.Range("B5").NumberFormat = "@"
$ws3.Range("B5").Value = "7.4"     # Avg LinkedIn Quality Score

$ws3.Range("B6").NumberFormat = "@"
$ws3.Range("B6").Value = "9.3"     # Avg Profile Completeness

Write-Host "Added Jacqueline Shen and Nikki Gittins; refreshed Summary metrics."
